$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the review "blue" flag for row 10 from "yes" to "no"
$ws.Range("G10").Value = "no"

# Move the active selection to G10 (matching the post-edit cursor position)
$ws.Range("G10").Select()
